$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the existing row 354 ("Femacal de La Calera"
# week of 2020-12-02), pushing the remaining weekly blocks down by two rows.
$ws.Rows("354:355").Insert()

# New week's data (2021-09-09 -> serial 44448), "Primera" / "Segunda" quality rows.
$ws.Range("A354").Value = 3
$ws.Range("B354").Value = "Femacal de La Calera"
$ws.Range("C354").Value = "Coquimbo"
$ws.Range("D354").Value = 44448
$ws.Range("E354").Value = 5
$ws.Range("F354").Value = 100112023
$ws.Range("G354").Value = "Brócoli"
$ws.Range("H354").Value = "Sin especificar"
$ws.Range("I354").Value = "Primera"
$ws.Range("J354").Value = 1600
$ws.Range("K354").Value = 600
$ws.Range("L354").Value = 600
$ws.Range("M354").Value = 600
$ws.Range("N354").Value = "$/unidad"
$ws.Range("O354").Value = "Provincia de Quillota"
$ws.Range("P354").Value = 600
$ws.Range("Q354").Value = 1
$ws.Range("R354").Value = "Hortaliza"

$ws.Range("A355").Value = 3
$ws.Range("B355").Value = "Femacal de La Calera"
$ws.Range("C355").Value = "Coquimbo"
$ws.Range("D355").Value = 44448
$ws.Range("E355").Value = 5
$ws.Range("F355").Value = 100112023
$ws.Range("G355").Value = "Brócoli"
$ws.Range("H355").Value = "Sin especificar"
$ws.Range("I355").Value = "Segunda"
$ws.Range("J355").Value = 2450
$ws.Range("K355").Value = 450
$ws.Range("L355").Value = 550
$ws.Range("M355").Value = 511
$ws.Range("N355").Value = "$/unidad"
$ws.Range("O355").Value = "Provincia de Quillota"
$ws.Range("P355").Value = 511
$ws.Range("Q355").Value = 1
$ws.Range("R355").Value = "Hortaliza"
